$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells being updated with new crypto price/volume data.
# NumberFormat is forced to Text ("@") before assignment so that
# plain-decimal-looking values (e.g. "302.95") are not silently
# coerced into numeric cells by Excel, matching the source data
# which stores every Price/Volume cell as text. The style is reset
# back to Normal afterwards so no extra formatting is introduced.
$cells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E9", "D10", "E10",
    "E11", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18",
    "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "E24", "E25", "E26", "D27",
    "E27", "D29", "E29", "E30", "D31", "E31", "E32", "E33", "D34", "E34", "E35", "E36", "D37",
    "E37", "E38", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45",
    "D46", "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)

foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.593.83'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '2.295.70'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '302.95'
$ws.Range('E5').Value = '  -2.55%  '
$ws.Range('D6').Value = '98.73'
$ws.Range('E6').Value = '  -6.69%  '
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -5.02%  '
$ws.Range('E9').Value = '  -5.65%  '
$ws.Range('D10').Value = '34.51'
$ws.Range('E10').Value = '  -6.81%  '
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  -4.29%  '
$ws.Range('D14').Value = '15.73'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').Value = '2.652.94'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '2.313.67'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '0.801'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = '42.525.67'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('D19').Value = '0.0₃0901'
$ws.Range('E19').Value = '  -3.27%  '
$ws.Range('D20').Value = '11.47'
$ws.Range('E20').Value = '  -6.18%  '
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').Value = '67.92'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').Value = '235.26'
$ws.Range('E23').Value = '  -3.17%  '
$ws.Range('E24').Value = '  -3.62%  '
$ws.Range('E25').Value = '  -4.05%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '25.00'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D29').Value = '34.79'
$ws.Range('E29').Value = '  -6.82%  '
$ws.Range('E30').Value = '  -5.26%  '
$ws.Range('D31').Value = '163.52'
$ws.Range('E31').Value = '  -1.75%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  -5.82%  '
$ws.Range('D34').Value = '4.62'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').Value = '  -4.97%  '
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('D37').Value = '16.94'
$ws.Range('E37').Value = '  -8.30%  '
$ws.Range('E38').Value = '  -6.26%  '
$ws.Range('E39').Value = '  -4.65%  '
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  -6.59%  '
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('D42').Value = '2.39'
$ws.Range('E42').Value = '  -12.14%  '
$ws.Range('D43').Value = '1.982.42'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('E44').Value = '  -4.99%  '
$ws.Range('D45').Value = '18.55'
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('D46').Value = '10.22'
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('E47').Value = '  -8.26%  '
$ws.Range('D48').Value = '55.35'
$ws.Range('E48').Value = '  -3.12%  '
$ws.Range('D49').Value = '2.85'
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('D50').Value = '2.521.16'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').Value = '4.68'
$ws.Range('E51').Value = '  -0.85%  '

foreach ($c in $cells) {
    $ws.Range($c).Style = "Normal"
}
